$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 236, shifting existing rows 236..313 down to 237..314.
$ws.Rows.Item(236).Insert()

# Populate the newly inserted row 236 with the new record.
$ws.Cells.Item(236, 1).Value = 10
$ws.Cells.Item(236, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(236, 3).Value = "La Araucanía"
$ws.Cells.Item(236, 4).Value = 45215
$ws.Cells.Item(236, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(236, 5).Value = 9
$ws.Cells.Item(236, 6).Value = 100112012
$ws.Cells.Item(236, 7).Value = "Espinaca"
$ws.Cells.Item(236, 8).Value = "Sin especificar"
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 55
$ws.Cells.Item(236, 11).Value = 12000
$ws.Cells.Item(236, 12).Value = 12000
$ws.Cells.Item(236, 13).Value = 12000
$ws.Cells.Item(236, 14).Value = "$/docena de atados"
$ws.Cells.Item(236, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(236, 16).Value = 4000
$ws.Cells.Item(236, 17).Value = 3
$ws.Cells.Item(236, 18).Value = "Hortaliza"
